$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1 (o_10): add new header/column, edit existing content ---
$s1_a2 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0 0
 G 0 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 1 0 0 1 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
Solution: A -> E -> F -> G -> H -> L -> P
        
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 0 0 0 1 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 0 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 1 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 1 0 1 0 0 1 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1 0
 L 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
    
"@
$s1_b2 = @"
A -> E -> I -> J -> K -> O -> P
"@
$s1_c2 = @"
The shortest path from node A to node P is A -> E -> F -> G -> H -> L -> P.
"@
$s1_d2 = @"
Wrong
"@
$s1_e2 = @"
Output: 2/6
"@

$ws1.Range("E1").Value = "evaluator_partial_correctness"
$headerRange1 = $ws1.Range("A1:E1")
$headerRange1.Font.Bold = $true
$headerRange1.HorizontalAlignment = -4108
$headerRange1.VerticalAlignment = -4160
$headerRange1.Borders.LineStyle = 1

$ws1.Range("A2").Value = $s1_a2
$ws1.Range("B2").Value = $s1_b2
$ws1.Range("C2").Value = $s1_c2
$ws1.Range("D2").Value = $s1_d2
$ws1.Range("E2").Value = $s1_e2

# --- sheet2 (o_20): new sheet ---
$s2_a2 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 24 nodes labelled A to X. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
Solution: A -> F -> G -> L -> M -> N -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node X?
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@
$s2_b2 = @"
A -> F -> J -> O -> T -> U -> V -> W -> X
"@
$s2_c2 = @"
The shortest path from node A to node X is A -> B -> C -> H -> R -> X.
"@
$s2_d2 = @"
Wrong
"@
$s2_e2 = @"
Output: 1/6
"@

$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "o_20"
$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$headerRange2 = $ws2.Range("A1:E1")
$headerRange2.Font.Bold = $true
$headerRange2.HorizontalAlignment = -4108
$headerRange2.VerticalAlignment = -4160
$headerRange2.Borders.LineStyle = 1

$ws2.Range("A2").Value = $s2_a2
$ws2.Range("B2").Value = $s2_b2
$ws2.Range("C2").Value = $s2_c2
$ws2.Range("D2").Value = $s2_d2
$ws2.Range("E2").Value = $s2_e2

# --- sheet3 (o_20_jumbled): new sheet ---
$s3_a2 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
Solution: A -> F -> G -> L -> M -> N -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@
$s3_b2 = @"
A -> F -> G -> L -> M -> N -> O -> T -> Y
"@
$s3_c2 = @"
To find the shortest path from node A to node Y, we can use the breadth-first search algorithm.
Starting from node A, we explore its neighboring nodes (B and F). We mark these nodes as visited and add them to a queue. 
Then, we repeat the process for the next level of nodes in the queue. For node B, we explore its neighboring nodes (A, C, and G). We mark the new nodes as visited and add them to the queue.
We continue this process until we reach node Y. Once we reach node Y, we can backtrack the shortest path by following the parent pointers from Y back to A.
Here is the step-by-step process:
1. Start at node A.
2. Add node A to the visited set and enqueue it.
3. While the queue is not empty:
   a. Dequeue a node from the queue.
   b. If the dequeued node is Y, we have found the shortest path. Break.
   c. For each neighbor of the dequeued node that is not visited:
      i. Mark the neighbor as visited.
      ii. Enqueue the neighbor.
      iii. Set the parent of the neighbor as the dequeued node.
4. Backtrack from node Y to node A using the parent pointers to determine the shortest path.
Based on the adjacency matrix provided, the shortest path from node A to node Y is: A -> F -> G -> L -> M -> N -> O -> T -> Y.
"@
$s3_d2 = @"
Correct
"@
$s3_e2 = @"
Output: 9/9
"@

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "o_20_jumbled"
$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$headerRange3 = $ws3.Range("A1:E1")
$headerRange3.Font.Bold = $true
$headerRange3.HorizontalAlignment = -4108
$headerRange3.VerticalAlignment = -4160
$headerRange3.Borders.LineStyle = 1

$ws3.Range("A2").Value = $s3_a2
$ws3.Range("B2").Value = $s3_b2
$ws3.Range("C2").Value = $s3_c2
$ws3.Range("D2").Value = $s3_d2
$ws3.Range("E2").Value = $s3_e2

$ws1.Select()
